$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "Price"
$ws.Range("B5").Value = "₹38,990"
